$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2028.6666
$ws.Range("I18").Value = 2028.6666
$ws.Range("K18").Value = 2028.6666
$ws.Range("M18").Value = -1744.6666
$ws.Range("H40").Value = 1276.25
$ws.Range("I40").Value = 1114.1111
$ws.Range("J40").Value = 1484.7142
$ws.Range("K40").Value = 1114.1111
$ws.Range("L40").Value = 1484.7142
$ws.Range("M40").Value = -939.1111000000001
$ws.Range("N40").Value = -1834.7142
$ws.Range("H43").Value = 3243.5557
$ws.Range("J43").Value = 3398.8333
$ws.Range("L43").Value = 3398.8333
$ws.Range("N43").Value = -3536.8333
$ws.Range("H112").Value = 3634.1428
$ws.Range("J112").Value = 4165.4443
$ws.Range("L112").Value = 12496.3329
$ws.Range("N112").Value = -14712.3329
$ws.Range("H114").Value = 80000
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H132").Value = 7103.543
$ws.Range("I132").Value = 8190.967
$ws.Range("K132").Value = 24572.901
$ws.Range("M132").Value = -22042.901
$ws.Range("H137").Value = 13162081
$ws.Range("I137").Value = 15153063
$ws.Range("K137").Value = 45459189
$ws.Range("M137").Value = -45456639
$ws.Range("H138").Value = 9132.286
$ws.Range("I138").Value = 9251.4
$ws.Range("J138").Value = 6750
$ws.Range("K138").Value = 27754.2
$ws.Range("L138").Value = 20250
$ws.Range("M138").Value = -22614.2
$ws.Range("N138").Value = -30530
$ws.Range("H141").Value = 1277.4286
$ws.Range("J141").Value = 1500
$ws.Range("L141").Value = 4500
$ws.Range("N141").Value = -14860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6672332.5
$ws.Range("I61").Value = 6328.6924
$ws.Range("J61").Value = 50001356
$ws.Range("K61").Value = 6328.6924
$ws.Range("L61").Value = 50001356
$ws.Range("M61").Value = -6116.6924
$ws.Range("N61").Value = -50001780
$ws.Range("H74").Value = 1165614.9
$ws.Range("I74").Value = 2062242.8
$ws.Range("K74").Value = 2062242.8
$ws.Range("M74").Value = -2061368.8
$ws.Range("H77").Value = 1165614.9
$ws.Range("I77").Value = 2062242.8
$ws.Range("K77").Value = 10311214
$ws.Range("M77").Value = -10306846
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("N107").Value = 0
$ws.Range("H113").Value = 73199
$ws.Range("J113").Value = 73199
$ws.Range("L113").Value = 73199
$ws.Range("N113").Value = -81877
$ws.Range("H132").Value = 783806.75
$ws.Range("I132").Value = 963880.4399999999
$ws.Range("J132").Value = 3487.5
$ws.Range("K132").Value = 2891641.32
$ws.Range("L132").Value = 10462.5
$ws.Range("M132").Value = -2889111.32
$ws.Range("N132").Value = -15522.5
$ws.Range("H136").Value = 6672332.5
$ws.Range("I136").Value = 6328.6924
$ws.Range("J136").Value = 50001356
$ws.Range("K136").Value = 18986.0772
$ws.Range("L136").Value = 150004068
$ws.Range("M136").Value = -16436.0772
$ws.Range("N136").Value = -150009168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 740.6667
$ws.Range("I107").Value = 740.6667
$ws.Range("K107").Value = 740.6667
$ws.Range("M107").Value = 1179.3333
$ws.Range("H134").Value = 8786506
$ws.Range("I134").Value = 10519.909
$ws.Range("K134").Value = 31559.727
$ws.Range("M134").Value = -29024.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 86.42856999999999
$ws.Range("I7").Value = 45.583332
$ws.Range("K7").Value = 45.583332
$ws.Range("M7").Value = 67.416668
$ws.Range("H58").Value = 6099954
$ws.Range("I58").Value = 5053659
$ws.Range("K58").Value = 5053659
$ws.Range("M58").Value = -5053456
$ws.Range("H136").Value = 6099954
$ws.Range("I136").Value = 5053659
$ws.Range("K136").Value = 15160977
$ws.Range("M136").Value = -15158427

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 100000130
$ws.Range("J11").Value = 250000240
$ws.Range("L11").Value = 750000720
$ws.Range("N11").Value = -750001000
$ws.Range("H23").Value = 453.73685
$ws.Range("I23").Value = 441.4
$ws.Range("J23").Value = 458.14285
$ws.Range("K23").Value = 1324.2
$ws.Range("L23").Value = 1374.42855
$ws.Range("M23").Value = -1089.2
$ws.Range("N23").Value = -1844.42855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.666664
$ws.Range("I2").Value = 47.6
$ws.Range("K2").Value = 47.6
$ws.Range("M2").Value = 65.40000000000001
$ws.Range("H70").Value = 57373.832
$ws.Range("I70").Value = 72250
$ws.Range("J70").Value = 42497.668
$ws.Range("K70").Value = 72250
$ws.Range("L70").Value = 42497.668
$ws.Range("M70").Value = -71980
$ws.Range("N70").Value = -43037.668
$ws.Range("H73").Value = 57373.832
$ws.Range("I73").Value = 72250
$ws.Range("J73").Value = 42497.668
$ws.Range("K73").Value = 72250
$ws.Range("L73").Value = 42497.668
$ws.Range("M73").Value = -71314
$ws.Range("N73").Value = -44369.668
$ws.Range("H141").Value = 93333
$ws.Range("J141").Value = 93333
$ws.Range("L141").Value = 93333
$ws.Range("N141").Value = -103693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H114").Value = 61851.668
$ws.Range("J114").Value = 61851.668
$ws.Range("L114").Value = 61851.668
$ws.Range("N114").Value = -70529.66800000001
$ws.Range("H132").Value = 2154057
$ws.Range("I132").Value = 2901959.2
$ws.Range("K132").Value = 8705877.600000001
$ws.Range("M132").Value = -8703347.600000001
$ws.Range("H136").Value = 13891024
$ws.Range("I136").Value = 8930825
$ws.Range("J136").Value = 31251722
$ws.Range("K136").Value = 26792475
$ws.Range("L136").Value = 93755166
$ws.Range("M136").Value = -26789925
$ws.Range("N136").Value = -93760266

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 125004424
$ws.Range("J81").Value = 250004350
$ws.Range("L81").Value = 500008700
$ws.Range("N81").Value = -500010822
$ws.Range("H84").Value = 125004424
$ws.Range("J84").Value = 250004350
$ws.Range("L84").Value = 2500043500
$ws.Range("N84").Value = -2500054108
$ws.Range("H132").Value = 5748951
$ws.Range("I132").Value = 7937886.5
$ws.Range("J132").Value = 2995.625
$ws.Range("K132").Value = 23813659.5
$ws.Range("L132").Value = 8986.875
$ws.Range("M132").Value = -23811129.5
$ws.Range("N132").Value = -14046.875
